$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 9007.714
$ws.Cells.Item(2, 9).Value = 1392.8889
$ws.Cells.Item(2, 10).Value = 22714.4
$ws.Cells.Item(2, 11).Value = 1392.8889
$ws.Cells.Item(2, 12).Value = 22714.4
$ws.Cells.Item(2, 13).Value = -1279.8889
$ws.Cells.Item(2, 14).Value = -22940.4
$ws.Cells.Item(33, 8).Value = 5556075.5
$ws.Cells.Item(33, 9).Value = 8333960.5
$ws.Cells.Item(33, 10).Value = 305.33334
$ws.Cells.Item(33, 11).Value = 8333960.5
$ws.Cells.Item(33, 12).Value = 305.33334
$ws.Cells.Item(33, 13).Value = -8333731.5
$ws.Cells.Item(33, 14).Value = -763.33334
$ws.Cells.Item(100, 8).Value = 65999.75
$ws.Cells.Item(100, 10).Value = 7999
$ws.Cells.Item(100, 12).Value = 7999
$ws.Cells.Item(100, 14).Value = -9081
$ws.Cells.Item(112, 8).Value = 2903.84
$ws.Cells.Item(112, 10).Value = 3046.182
$ws.Cells.Item(112, 12).Value = 9138.545999999998
$ws.Cells.Item(112, 14).Value = -11354.546
$ws.Cells.Item(113, 8).Value = 2877.4443
$ws.Cells.Item(113, 10).Value = 2943.6667
$ws.Cells.Item(113, 12).Value = 2943.6667
$ws.Cells.Item(113, 14).Value = -9451.6667
$ws.Cells.Item(115, 8).Value = 893.4
$ws.Cells.Item(115, 9).Value = 893.4
$ws.Cells.Item(115, 11).Value = 2680.2
$ws.Cells.Item(115, 13).Value = -1113.2
$ws.Cells.Item(137, 8).Value = 13896205
$ws.Cells.Item(137, 9).Value = 25001146
$ws.Cells.Item(137, 10).Value = 15028.125
$ws.Cells.Item(137, 11).Value = 75003438
$ws.Cells.Item(137, 12).Value = 45084.375
$ws.Cells.Item(137, 13).Value = -75000888
$ws.Cells.Item(137, 14).Value = -50184.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 18184500
$ws.Cells.Item(2, 10).Value = 1010
$ws.Cells.Item(2, 12).Value = 1010
$ws.Cells.Item(2, 14).Value = -1236
$ws.Cells.Item(32, 8).Value = 49313.957
$ws.Cells.Item(32, 9).Value = 49313.957
$ws.Cells.Item(32, 11).Value = 49313.957
$ws.Cells.Item(32, 13).Value = -49026.957
$ws.Cells.Item(80, 8).Value = 87554.25
$ws.Cells.Item(80, 9).Value = 20000
$ws.Cells.Item(80, 10).Value = 110072.336
$ws.Cells.Item(80, 11).Value = 20000
$ws.Cells.Item(80, 12).Value = 110072.336
$ws.Cells.Item(80, 13).Value = -19002
$ws.Cells.Item(80, 14).Value = -112068.336
$ws.Cells.Item(83, 8).Value = 87554.25
$ws.Cells.Item(83, 9).Value = 20000
$ws.Cells.Item(83, 10).Value = 110072.336
$ws.Cells.Item(83, 11).Value = 60000
$ws.Cells.Item(83, 12).Value = 330217.008
$ws.Cells.Item(83, 13).Value = -55008
$ws.Cells.Item(83, 14).Value = -340201.008
$ws.Cells.Item(110, 8).Value = 128336340
$ws.Cells.Item(110, 9).Value = 128336340
$ws.Cells.Item(110, 11).Value = 128336340
$ws.Cells.Item(110, 13).Value = -128334295
$ws.Cells.Item(116, 8).Value = 18184500
$ws.Cells.Item(116, 10).Value = 1010
$ws.Cells.Item(116, 12).Value = 1010
$ws.Cells.Item(116, 14).Value = -5598
$ws.Cells.Item(122, 8).Value = 1664.95
$ws.Cells.Item(122, 9).Value = 1605.8823
$ws.Cells.Item(122, 11).Value = 4817.6469
$ws.Cells.Item(122, 13).Value = -2367.6469
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 18184500
$ws.Cells.Item(3, 10).Value = 1010
$ws.Cells.Item(3, 12).Value = 1010
$ws.Cells.Item(3, 14).Value = -1238
$ws.Cells.Item(102, 8).Value = 12674.25
$ws.Cells.Item(102, 10).Value = 42612
$ws.Cells.Item(102, 12).Value = 42612
$ws.Cells.Item(102, 14).Value = -49102
$ws.Cells.Item(116, 8).Value = 742
$ws.Cells.Item(116, 10).Value = 742
$ws.Cells.Item(116, 12).Value = 742
$ws.Cells.Item(116, 14).Value = -9920
$ws.Cells.Item(137, 8).Value = 126663.5
$ws.Cells.Item(137, 10).Value = 126663.5
$ws.Cells.Item(137, 12).Value = 126663.5
$ws.Cells.Item(137, 14).Value = -136863.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 58828520
$ws.Cells.Item(31, 9).Value = 166667500
$ws.Cells.Item(31, 11).Value = 166667500
$ws.Cells.Item(31, 13).Value = -166667205
$ws.Cells.Item(34, 8).Value = 58828520
$ws.Cells.Item(34, 9).Value = 166667500
$ws.Cells.Item(34, 11).Value = 166667500
$ws.Cells.Item(34, 13).Value = -166667298
$ws.Cells.Item(99, 8).Value = 4879.6665
$ws.Cells.Item(99, 9).Value = 3323.75
$ws.Cells.Item(99, 11).Value = 3323.75
$ws.Cells.Item(99, 13).Value = -1825.75
$ws.Cells.Item(126, 8).Value = 4879.6665
$ws.Cells.Item(126, 9).Value = 3323.75
$ws.Cells.Item(126, 11).Value = 9971.25
$ws.Cells.Item(126, 13).Value = -7501.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(111, 8).Value = 16018
$ws.Cells.Item(113, 8).Value = 672.5
$ws.Cells.Item(113, 9).Value = 299.25
$ws.Cells.Item(113, 11).Value = 897.75
$ws.Cells.Item(113, 13).Value = 1272.25
$ws.Cells.Item(122, 8).Value = 13486.833
$ws.Cells.Item(122, 9).Value = 35710
$ws.Cells.Item(122, 11).Value = 321390
$ws.Cells.Item(122, 13).Value = -318940
$ws.Cells.Item(129, 8).Value = 1277.9333
$ws.Cells.Item(129, 9).Value = 654.9286
$ws.Cells.Item(129, 11).Value = 1964.7858
$ws.Cells.Item(129, 13).Value = 3035.2142
$ws.Cells.Item(131, 8).Value = 7039.231
$ws.Cells.Item(131, 9).Value = 963.6667
$ws.Cells.Item(131, 10).Value = 8861.9
$ws.Cells.Item(131, 11).Value = 2891.0001
$ws.Cells.Item(131, 12).Value = 26585.7
$ws.Cells.Item(131, 13).Value = 2148.9999
$ws.Cells.Item(131, 14).Value = -36665.7
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 53541.05
$ws.Cells.Item(2, 9).Value = 1224.7142
$ws.Cells.Item(2, 11).Value = 1224.7142
$ws.Cells.Item(2, 13).Value = -1111.7142
$ws.Cells.Item(46, 8).Value = 33395
$ws.Cells.Item(46, 9).Value = 15041
$ws.Cells.Item(46, 10).Value = 37983.5
$ws.Cells.Item(46, 11).Value = 15041
$ws.Cells.Item(46, 12).Value = 37983.5
$ws.Cells.Item(46, 13).Value = -14885
$ws.Cells.Item(46, 14).Value = -38295.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(5, 8).Value = 14899
$ws.Cells.Item(5, 10).Value = 14899
$ws.Cells.Item(5, 12).Value = 14899
$ws.Cells.Item(5, 14).Value = -15125
$ws.Cells.Item(40, 8).Value = 37504
$ws.Cells.Item(40, 9).Value = 37504
$ws.Cells.Item(40, 11).Value = 37504
$ws.Cells.Item(40, 13).Value = -37368
$ws.Cells.Item(46, 8).Value = 4095.02
$ws.Cells.Item(46, 9).Value = 1476.6666
$ws.Cells.Item(46, 10).Value = 5217.1714
$ws.Cells.Item(46, 11).Value = 1476.6666
$ws.Cells.Item(46, 12).Value = 5217.1714
$ws.Cells.Item(46, 13).Value = -1288.6666
$ws.Cells.Item(46, 14).Value = -5593.1714
$ws.Cells.Item(55, 8).Value = 182.28572
$ws.Cells.Item(55, 9).Value = 135.29411
$ws.Cells.Item(55, 11).Value = 135.29411
$ws.Cells.Item(55, 13).Value = 37.70589000000001
$ws.Cells.Item(68, 8).Value = 2585.2
$ws.Cells.Item(68, 9).Value = 2585.2
$ws.Cells.Item(68, 11).Value = 2585.2
$ws.Cells.Item(68, 13).Value = -1836.2
$ws.Cells.Item(71, 8).Value = 2585.2
$ws.Cells.Item(71, 9).Value = 2585.2
$ws.Cells.Item(71, 11).Value = 12926
$ws.Cells.Item(71, 13).Value = -9182
$ws.Cells.Item(76, 8).Value = 15999.333
$ws.Cells.Item(76, 10).Value = 15999.333
$ws.Cells.Item(76, 12).Value = 15999.333
$ws.Cells.Item(76, 14).Value = -16675.333
$ws.Cells.Item(79, 8).Value = 15999.333
$ws.Cells.Item(79, 10).Value = 15999.333
$ws.Cells.Item(79, 12).Value = 15999.333
$ws.Cells.Item(79, 14).Value = -18339.333
$ws.Cells.Item(124, 8).Value = 0
$ws.Cells.Item(124, 10).Value = 0
$ws.Cells.Item(124, 12).Value = 0
$ws.Cells.Item(124, 14).ClearContents()
$ws.Cells.Item(125, 8).Value = 183332.47
$ws.Cells.Item(125, 10).Value = 183332.47
$ws.Cells.Item(125, 12).Value = 183332.47
$ws.Cells.Item(125, 14).Value = -193172.47
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 103038.57
$ws.Cells.Item(124, 10).Value = 103038.57
$ws.Cells.Item(124, 12).Value = 103038.57
$ws.Cells.Item(124, 14).Value = -112858.57
